$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2  = @{ C = 224473.115141326;  D = 244036.25347797;   E = 237264.763743061; F = 211681.466539592; G = 204909.976804683 }
    3  = @{ C = 220814.913719902;  D = 243046.137472794;  E = 235351.129604082; F = 206278.697835723; G = 198583.68996701 }
    4  = @{ C = 257852.788003249;  D = 286886.682543403;  E = 276837.031231842; F = 238868.544774655; G = 228818.893463094 }
    5  = @{ C = 238543.07305751;   D = 267974.440705653;  E = 257787.20998202;  F = 219298.936133;    G = 209111.705409367 }
    6  = @{ C = 275565.933415601;  D = 312298.259488872;  E = 299583.910458349; F = 251547.956372854; G = 238833.607342331 }
    7  = @{ C = 278613.164291491;  D = 318324.806243271;  E = 304579.211277601; F = 252647.117305382; G = 238901.522339712 }
    8  = @{ C = 278950.576958915;  D = 321131.043705324;  E = 306530.901710111; F = 251370.252207719; G = 236770.110212506 }
    9  = @{ C = 244575.86980338;   D = 283566.848241238;  E = 270070.700245828; F = 219081.039360933; G = 205584.891365522 }
    10 = @{ C = 214225.214332171;  D = 250051.536559116;  E = 237650.787301779; F = 190799.641362563; G = 178398.892105226 }
    11 = @{ C = 213620.918937558;  D = 250941.835969952;  E = 238023.754922143; F = 189218.082952973; G = 176300.001905164 }
    12 = @{ C = 223330.017759523;  D = 263947.932313949;  E = 249888.644549085; F = 196771.390969961; G = 182712.103205097 }
    13 = @{ C = 225313.159432494;  D = 267846.457004606;  E = 253124.187860195; F = 197502.131004793; G = 182779.861860383 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
